# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-24 11:13:48
# Reorder "Recorded By" email lists, update recorded/missing attendance counts,
# and refresh the resulting "Average Attendance %" figures on the
# "Session Analysis Results" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - Recorded By
$ws.Range("G2").Value = "System, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 3 - Recorded By
$ws.Range("G3").Value = "System, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 - Recorded By
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - Recorded By
$ws.Range("G5").Value = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 7 - Recorded By + Students count
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("H7").Value = "41/251"

# Row 10 - updated average attendance percentage (stored as literal text,
# so format it as Text first to keep Excel from re-interpreting the "%"
# string as a numeric percentage)
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "25.2%"

# Row 15 - Recorded By + average attendance percentage
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "25.2%"

$wb.Save()
